# Scheduled-runner refresh of per-leve market price/profit figures
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 118749.25
$ws.Range("J120").Value = 118749.25
$ws.Range("L120").Value = 118749.25
$ws.Range("N120").Value = -128425.25

$ws.Range("H129").Value = 706.7826
$ws.Range("J129").Value = 1023.5833
$ws.Range("L129").Value = 3070.7499
$ws.Range("N129").Value = -13070.7499

$ws.Range("H135").Value = 965.8
$ws.Range("I135").Value = 369.7451
$ws.Range("J135").Value = 4343.4443
$ws.Range("K135").Value = 3327.7059
$ws.Range("L135").Value = 39090.9987
$ws.Range("M135").Value = -792.7058999999999
$ws.Range("N135").Value = -44160.9987

$ws.Range("H136").Value = 52566.152
$ws.Range("J136").Value = 52566.152
$ws.Range("L136").Value = 52566.152
$ws.Range("N136").Value = -62766.152

$ws.Range("H138").Value = 1888.9302
$ws.Range("I138").Value = 862.07465
$ws.Range("J138").Value = 5509.9473
$ws.Range("K138").Value = 2586.22395
$ws.Range("L138").Value = 16529.8419
$ws.Range("M138").Value = 2553.77605
$ws.Range("N138").Value = -26809.8419

$ws.Range("H139").Value = 74915.336
$ws.Range("J139").Value = 74915.336
$ws.Range("L139").Value = 74915.336
$ws.Range("N139").Value = -85195.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3655.45
$ws.Range("I32").Value = 2365.5889
$ws.Range("J32").Value = 15264.2
$ws.Range("K32").Value = 2365.5889
$ws.Range("L32").Value = 15264.2
$ws.Range("M32").Value = -2078.5889
$ws.Range("N32").Value = -15838.2

$ws.Range("H52").Value = 44890
$ws.Range("J52").Value = 44890
$ws.Range("L52").Value = 44890
$ws.Range("N52").Value = -45526

$ws.Range("H61").Value = 1022.43396
$ws.Range("I61").Value = 775.9761999999999
$ws.Range("J61").Value = 1963.4546
$ws.Range("K61").Value = 775.9761999999999
$ws.Range("L61").Value = 1963.4546
$ws.Range("M61").Value = -563.9761999999999
$ws.Range("N61").Value = -2387.4546

$ws.Range("H132").Value = 2339.4722
$ws.Range("I132").Value = 2458.6
$ws.Range("K132").Value = 7375.799999999999
$ws.Range("M132").Value = -4845.799999999999

$ws.Range("H136").Value = 1022.43396
$ws.Range("I136").Value = 775.9761999999999
$ws.Range("J136").Value = 1963.4546
$ws.Range("K136").Value = 2327.9286
$ws.Range("L136").Value = 5890.3638
$ws.Range("M136").Value = 222.0714000000003
$ws.Range("N136").Value = -10990.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1513.9231
$ws.Range("I107").Value = 1204.75
$ws.Range("K107").Value = 1204.75
$ws.Range("M107").Value = 715.25

$ws.Range("H132").Value = 45122.5
$ws.Range("J132").Value = 45122.5
$ws.Range("L132").Value = 45122.5
$ws.Range("N132").Value = -55242.5

$ws.Range("H134").Value = 989.3226
$ws.Range("I134").Value = 851.7083
$ws.Range("K134").Value = 2555.1249
$ws.Range("M134").Value = -20.1248999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 3706.5
$ws.Range("I39").Value = 1850
$ws.Range("J39").Value = 4634.75
$ws.Range("K39").Value = 1850
$ws.Range("L39").Value = 4634.75
$ws.Range("M39").Value = -1459
$ws.Range("N39").Value = -5416.75

$ws.Range("H49").Value = 3706.5
$ws.Range("I49").Value = 1850
$ws.Range("J49").Value = 4634.75
$ws.Range("K49").Value = 1850
$ws.Range("L49").Value = 4634.75
$ws.Range("M49").Value = -1668
$ws.Range("N49").Value = -4998.75

$ws.Range("H58").Value = 1013.9178
$ws.Range("I58").Value = 949.80646
$ws.Range("J58").Value = 1375.2727
$ws.Range("K58").Value = 949.80646
$ws.Range("L58").Value = 1375.2727
$ws.Range("M58").Value = -746.80646
$ws.Range("N58").Value = -1781.2727

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H122").Value = 3008.375
$ws.Range("I122").Value = 3256.4614
$ws.Range("J122").Value = 1933.3334
$ws.Range("K122").Value = 9769.3842
$ws.Range("L122").Value = 5800.0002
$ws.Range("M122").Value = -7319.3842
$ws.Range("N122").Value = -10700.0002

$ws.Range("H134").Value = 890.0328
$ws.Range("I134").Value = 890.0328
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2670.0984
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -135.0983999999999
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 1013.9178
$ws.Range("I136").Value = 949.80646
$ws.Range("J136").Value = 1375.2727
$ws.Range("K136").Value = 2849.41938
$ws.Range("L136").Value = 4125.8181
$ws.Range("M136").Value = -299.4193800000003
$ws.Range("N136").Value = -9225.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 2018
$ws.Range("I102").Value = 2018
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 6054
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3620
$ws.Range("N102").ClearContents()

$ws.Range("H118").Value = 1848.4
$ws.Range("I118").Value = 550
$ws.Range("J118").Value = 2578.75
$ws.Range("K118").Value = 1650
$ws.Range("L118").Value = 7736.25
$ws.Range("M118").Value = -407
$ws.Range("N118").Value = -10222.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1130.65
$ws.Range("I126").Value = 850.8182
$ws.Range("J126").Value = 1472.6666
$ws.Range("K126").Value = 2552.4546
$ws.Range("L126").Value = 4417.9998
$ws.Range("M126").Value = -82.45460000000003
$ws.Range("N126").Value = -9357.9998

$ws.Range("H132").Value = 1527.9143
$ws.Range("I132").Value = 1644.0344
$ws.Range("J132").Value = 966.6667
$ws.Range("K132").Value = 4932.1032
$ws.Range("L132").Value = 2900.0001
$ws.Range("M132").Value = -2402.1032
$ws.Range("N132").Value = -7960.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1968
$ws.Range("I16").Value = 1027.4546
$ws.Range("J16").Value = 5416.6665
$ws.Range("K16").Value = 1027.4546
$ws.Range("L16").Value = 5416.6665
$ws.Range("M16").Value = -857.4546
$ws.Range("N16").Value = -5756.6665

$ws.Range("H46").Value = 798.8333
$ws.Range("I46").Value = 723.5
$ws.Range("J46").Value = 949.5
$ws.Range("K46").Value = 723.5
$ws.Range("L46").Value = 949.5
$ws.Range("M46").Value = -535.5
$ws.Range("N46").Value = -1325.5

$ws.Range("H132").Value = 2828.2727
$ws.Range("I132").Value = 2430.4814
$ws.Range("K132").Value = 7291.4442
$ws.Range("M132").Value = -4761.4442

$ws.Range("H136").Value = 2195.9275
$ws.Range("I136").Value = 1533.3673
$ws.Range("K136").Value = 4600.1019
$ws.Range("M136").Value = -2050.1019

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 413.80264
$ws.Range("I132").Value = 363.19644
$ws.Range("J132").Value = 555.5
$ws.Range("K132").Value = 1089.58932
$ws.Range("L132").Value = 1666.5
$ws.Range("M132").Value = 1440.41068
$ws.Range("N132").Value = -6726.5
